$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '22.442.91'
$ws.Range("E2").Value = '  +8.95%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.609.68'
$ws.Range("E3").Value = '  +8.95%  '
$ws.Range("E4").Value = '  -0.71%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '304.55'
$ws.Range("E5").Value = '  +8.69%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9918'
$ws.Range("E6").Value = '  +3.87%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3694'
$ws.Range("E7").Value = '  +1.05%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3415'
$ws.Range("E8").Value = '  +11.44%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '42.33'
$ws.Range("E9").Value = '  +5.99%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.138'
$ws.Range("E10").Value = '  +7.20%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07087'
$ws.Range("E11").Value = '  +5.95%  '
$ws.Range("E12").Value = '  -0.30%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '19.80'
$ws.Range("E13").Value = '  +9.33%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.936'
$ws.Range("E14").Value = '  +7.41%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.667'
$ws.Range("E15").Value = '  +7.07%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.00001087'
$ws.Range("E16").Value = '  +4.99%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.603.31'
$ws.Range("E17").Value = '  +8.64%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.9924'
$ws.Range("E18").Value = '  +3.89%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06798'
$ws.Range("E19").Value = '  +14.09%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '78.39'
$ws.Range("E20").Value = '  +11.94%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.054'
$ws.Range("E21").Value = '  +9.95%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '16.11'
$ws.Range("E22").Value = '  +11.28%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '11.88'
$ws.Range("E23").Value = '  +7.16%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '22.459.97'
$ws.Range("E24").Value = '  +8.90%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.397'
$ws.Range("E25").Value = '  +5.90%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.553'
$ws.Range("E26").Value = '  +20.42%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '150.56'
$ws.Range("E27").Value = '  +5.07%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '19.60'
$ws.Range("E28").Value = '  +13.50%  '
$ws.Range("E29").Value = '  +9.02%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '122.32'
$ws.Range("E30").Value = '  +7.23%  '
$ws.Range("B31").Value = 'Filecoin'
$ws.Range("C31").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.134'
$ws.Range("E31").Value = '  +22.15%  '
$ws.Range("B32").Value = 'HuobiToken'
$ws.Range("C32").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.977'
$ws.Range("E32").Value = '  +0.24%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.9527'
$ws.Range("E33").Value = '  +17.15%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.08284'
$ws.Range("E34").Value = '  +4.00%  '
$ws.Range("E35").Value = '  +8.63%  '
$ws.Range("E36").Value = '  +15.50%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.266'
$ws.Range("E37").Value = '  +10.94%  '
$ws.Range("E38").Value = '  +3.64%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '8.631'
$ws.Range("E39").Value = '  +15.77%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.06101'
$ws.Range("E40").Value = '  +4.58%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.02230'
$ws.Range("E41").Value = '  +8.61%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.2030'
$ws.Range("E42").Value = '  +8.03%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.9921'
$ws.Range("E43").Value = '  +3.82%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.5937'
$ws.Range("E44").Value = '  +11.73%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.826'
$ws.Range("E45").Value = '  +8.03%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '13.23'
$ws.Range("E46").Value = '  +7.99%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5719'
$ws.Range("E47").Value = '  +10.07%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '127.50'
$ws.Range("E48").Value = '  +7.99%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.985'
$ws.Range("E49").Value = '  +8.74%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.06814'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '74.09'
$ws.Range("E51").Value = '  +9.49%  '
